$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 411, shifting the existing rows 411:475
# down to 412:476 (Excel copies formatting from the row above, which is
# what the target file shows for the date-formatted column D).
$ws.Rows.Item(411).Insert()

# The row that used to be 411 is now 412; duplicate its (unchanged)
# values into the freshly inserted row 411 first ...
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(411, $c).Value = $ws.Cells.Item(412, $c).Value2
}

# ... then overwrite the columns that actually carry new data for the
# new week's record (Fecha / Volumen / Precio minimo / Precio maximo /
# Precio promedio ponderado / Precio $/Kg).
$ws.Range("D411").Value = 45180
$ws.Range("J411").Value = 1160
$ws.Range("K411").Value = 450
$ws.Range("L411").Value = 500
$ws.Range("M411").Value = 475
$ws.Range("P411").Value = 950
